# xls2xml - completed series.json
# Update the "series" sheet: rename id_season header, add English
# translations to synopsis/season-name strings, and fix season 2's
# series id to match season 1's.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("series")
$ws.Activate()

# Header row: "id_season" -> "id season"
$ws.Range("D1").Value = "id season"

# Row 2 (season 1 of "A Lista" / "The List")
$ws.Range("C2").Value = "por:A melhor série de contagem regressiva de esportes. Com temas icônicos focados nos momentos e indivíduos mais memoráveis do esporte.|eng:English Synopsis"
$ws.Range("F2").Value = "por:Primeira temporada|eng:First season"

# Row 3 (season 2 of "A Lista" / "The List")
$ws.Range("A3").Value = "3d0666d2-0d6e-4687-b37b-1f65e173f889"
$ws.Range("C3").Value = "por:A melhor série de contagem regressiva de esportes. Com temas icônicos focados nos momentos e indivíduos mais memoráveis do esporte.|eng:English Synopsis"
$ws.Range("F3").Value = "por:Segunda temporada|eng:Second Season"

# Move the selection/active cell to F3
$ws.Range("F3").Select() | Out-Null
